$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the existing row 19 ("Baseline_2010_c31" row),
# pushing that row and everything below it down by two. Row 18 keeps its
# position but gains model/scenario labels; the first new row (19) becomes
# the "C109 without springs" data row, and the second new row (20) is left
# as a blank separator (matching the blank row that used to sit at the end
# of the "current" block).
$ws.Range("A19:A20").EntireRow.Insert()

# Row 18: label it as a CW3M / "Baseline 2010 C109" run, and give its
# mass-balance-discrepancy-fraction cell (R18) the same number format used
# by the other data rows instead of the separator-row format.
$ws.Range("A18").Value = "CW3M"
$ws.Range("B18").Value = "Baseline 2010 C109"
$ws.Range("R18").NumberFormat = "0.000000"

# Row 19: new "C109 without springs" data row.
$ws.Range("B19").Value = "C109 without springs"
$ws.Range("C19").Value = 2010
$ws.Range("D19").Value = 1090.199341
$ws.Range("E19").Value = 1990.4676509999999
$ws.Range("F19").Value = 1.2021059999999999
$ws.Range("G19").Value = 0
$ws.Range("H19").Value = 10.610913999999999
$ws.Range("I19").Value = 4.714251
$ws.Range("J19").Value = 8.8404570000000007
$ws.Range("K19").Value = 677.32849099999999
$ws.Range("L19").Value = 93.229797000000005
$ws.Range("M19").Value = 1114.0545649999999
$ws.Range("N19").Value = 1206.7479249999999
$ws.Range("O19").Value = 6798.6591799999997
$ws.Range("P19").Value = 29450.638672000001
$ws.Range("Q19").Value = 3.0069729999999999
$ws.Range("R19").Value = 0.00097099999999999997
$ws.Range("R19").NumberFormat = "0.000000"
$ws.Range("S19").Value = 2010

$ws.Range("B19").Select()
